# "Generate Report for Handoff"
#
# The report-generation step refreshed the handoff timestamp for the
# 28c8a306-8aa6-4b74-9a4e-4e62df69cd7f source file (a new handoff round
# was recorded, 23 seconds after the previous one) on all three sheets
# that surface it: the Overview summary sheet and the per-locale detail
# sheets (zh-cn, de-de). No other cell content changes.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-37-21 04:37:31"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-21 04:37:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-21 04:37:31"
